# Loan Provisioning 10 test cases
#
# 1) Swap the names of the two "Acc_Upfront" sheets that sit in slots
#    sheetId=9 (rId7) and sheetId=11 (rId9): what used to be called
#    "Acc_Upfront1" becomes "Acc_Upfront3", and what used to be called
#    "Acc_Upfront3" becomes "Acc_Upfront1". Tab order / underlying
#    sheetId/rId stay untouched - only the display names trade places.
# 2) Update the remembered cell selection on the "Transactions" sheet
#    from C6 to D2.
# 3) Update the remembered cell selection on the sheet that is now named
#    "Acc_Upfront1" (the tab that was/ is active) from E3 to I18.
# 4) Nudge the workbook's "first visible tab" scroll position forward by
#    one (best effort - harmless if the host doesn't track it).

$wb = $excel.ActiveWorkbook

# --- 1. Swap sheet names -------------------------------------------------
$wsWasUpfront1 = $wb.Worksheets.Item("Acc_Upfront1")
$wsWasUpfront3 = $wb.Worksheets.Item("Acc_Upfront3")

# Use a scratch name so the two swap without colliding.
$wsWasUpfront1.Name = "Acc_Upfront_swap_tmp"
$wsWasUpfront3.Name = "Acc_Upfront1"
$wsWasUpfront1.Name = "Acc_Upfront3"

# --- 2. Transactions sheet selection: C6 -> D2 ---------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Activate()
$wsTransactions.Range("D2").Select()

# --- 3. Active Acc_Upfront1 sheet selection: E3 -> I18 -------------------
$wsUpfront1 = $wb.Worksheets.Item("Acc_Upfront1")
$wsUpfront1.Activate()
$wsUpfront1.Range("I18").Select()

# --- 4. Scroll the visible tab strip forward by one (best effort) -------
try {
    $excel.ActiveWindow.DisplayedFirstSheet = 4
} catch {
    # Not fatal if the host doesn't expose/track this.
}

Write-Host "Sheets renamed, selections updated."
